$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "92.329.89", "0.430").
# Force text format first so Excel doesn't auto-convert numeric-looking
# values (and drop significant trailing zeros / thousands dots).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "92.329.89"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "3.095.48"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "236.58"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "610.76"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  -3.63%  "
$ws.Range("D8").Value = "0.389"
$ws.Range("E8").Value = "  +4.27%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "3.092.56"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").Value = "0.736"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "92.274.51"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "34.07"
$ws.Range("E15").Value = "  -3.50%  "
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "3.682.84"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "3.094.87"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "5.75"
$ws.Range("E21").Value = "  -4.51%  "
$ws.Range("D22").Value = "9.32"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "444.63"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").Value = "5.71"
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("D26").Value = "86.14"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").Value = "11.62"
$ws.Range("E27").Value = "  -3.68%  "
$ws.Range("D28").Value = "3.258.81"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "0.229"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").Value = "9.12"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "7.91"
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("D35").Value = "0.156"
$ws.Range("E35").Value = "  -8.08%  "
$ws.Range("D36").Value = "25.94"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("B37").Value = "MantraDAO"
$ws.Range("C37").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D37").Value = "3.89"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "486.08"
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("E41").Value = "  +7.80%  "
$ws.Range("D42").Value = "0.430"
$ws.Range("E42").Value = "  -5.23%  "
$ws.Range("D43").Value = "3.28"
$ws.Range("E43").Value = "  -5.01%  "
$ws.Range("D44").Value = "0.752"
$ws.Range("E44").Value = "  -24.75%  "
$ws.Range("D46").Value = "163.14"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").Value = "1.88"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").Value = "0.684"
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("E51").Value = "  -0.27%  "
